# App Clientes y correcciones en clases
#
# Fixes several "Reserva" (Si/No) values, renames a vehicle class, fixes a
# capitalization typo, and appends a brand-new booking row (Id 8, UTV ATV).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: writes $text into $range as a genuine text value, even when the
# text looks like a number or a date (e.g. "8", "2000", "12/12/2023") which
# Excel would otherwise auto-convert to a numeric/date cell. We do this by
# putting a text-returning formula ("=""...""") in a scratch cell far outside
# the used range, copying it, and pasting *values only* into the target -
# Excel keeps the pasted value as text instead of re-parsing it, and (unlike
# forcing a Text NumberFormat) this does not touch any cell styles.
function Set-TextValue($range, $text) {
    $scratch = $ws.Range("ZZ1")
    $escaped = $text -replace '"', '""'
    $scratch.Formula = '="' + $escaped + '"'
    $scratch.Copy()
    $range.PasteSpecial(-4163)
    $scratch.ClearContents()
}

# --- Corrections to existing rows ---
$ws.Range("L2").Value = "No"          # Reserva: Si -> No
$ws.Range("L3").Value = "Sí"          # Reserva: No -> Si
$ws.Range("B6").Value = "Naked Moto"  # Categoria Y Tipo: SUV Automovil -> Naked Moto
$ws.Range("L6").Value = "No"          # Reserva: Si -> No
$ws.Range("K7").Value = "No"          # Contraseña del conductor: NO -> No
$ws.Range("L8").Value = "Sí"          # Reserva: No -> Si

# --- New row: Id 8, UTV ATV ---
Set-TextValue $ws.Range("A9") "8"
$ws.Range("B9").Value = "UTV ATV"
Set-TextValue $ws.Range("C9") "12/12/2023"
$ws.Range("D9").Value = "Ruedas"
$ws.Range("E9").Value = "Ruedas"
$ws.Range("F9").Value = "13/12/2023"
$ws.Range("G9").Value = "Juan"
$ws.Range("H9").Value = "Cli"
Set-TextValue $ws.Range("I9") "0"
$ws.Range("J9").Value = "No"
$ws.Range("K9").Value = "No"
$ws.Range("L9").Value = "Sí"
Set-TextValue $ws.Range("M9") "2000"
